# ---------------------------------------------------------------------------
# Add the "2022-Q1" sheet (positioned between "2021-Q4" and "总计"), fill it
# with the fund-holdings detail table, and insert the summary row for the new
# quarter into the "总计" (Total) sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# Helper: write $val into $addr as literal TEXT (never let Excel's numeric
# auto-detection turn strings like "002345" or "30.12" into numbers / strip
# leading zeros). We stage the quoted literal in a scratch cell, then carry
# it over with a values-only paste so the destination cell keeps the default
# (unstyled) cell format -- exactly like the plain inlineStr cells already
# used on the neighbouring quarter sheets.
function Set-TextValue($sheet, $addr, $val) {
    $scratch = $sheet.Range("ZZ1")
    $scratch.Value = "'" + $val
    $scratch.Copy()
    $sheet.Range($addr).PasteSpecial(-4163)
    $scratch.Clear()
}

# --- create & position the new sheet ---------------------------------------
# NB: fetch worksheet references *after* Add() -- a reference grabbed before
# the insertion can end up pointing at the wrong tab once the collection
# shifts.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q1"
$q4 = $wb.Worksheets.Item("2021-Q4")
$newSheet.Move($null, $q4)

$ws = $wb.Worksheets.Item("2022-Q1")
$q4 = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")

# --- header row (B1:H1) + index column (A2:A22) formatting ------------------
# "2021-Q4" already has the exact look we need (bold/bordered/centred style)
# on its header row and its "A" index column, so copy formats from there.
$q4.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2").Copy()
$ws.Range("A2:A22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- header labels -----------------------------------------------------------
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# --- data rows ---------------------------------------------------------------
$ws.Range("A2").Value = 0
Set-TextValue $ws "B2" "160311"
Set-TextValue $ws "C2" "华夏蓝筹混合(LOF)"
Set-TextValue $ws "D2" "30.12"
Set-TextValue $ws "E2" "87.73"
Set-TextValue $ws "F2" "3.78"
Set-TextValue $ws "G2" "1.1385"
$ws.Range("H2").Value = 7

$ws.Range("A3").Value = 1
Set-TextValue $ws "B3" "400032"
Set-TextValue $ws "C3" "东方主题精选混合"
Set-TextValue $ws "D3" "29.31"
Set-TextValue $ws "E3" "85.31"
Set-TextValue $ws "F3" "3.06"
Set-TextValue $ws "G3" "0.8969"
$ws.Range("H3").Value = 7

$ws.Range("A4").Value = 2
Set-TextValue $ws "B4" "002345"
Set-TextValue $ws "C4" "华夏高端制造灵活配置混合"
Set-TextValue $ws "D4" "23.47"
Set-TextValue $ws "E4" "93.20"
Set-TextValue $ws "F4" "3.54"
Set-TextValue $ws "G4" "0.8308"
$ws.Range("H4").Value = 9

$ws.Range("A5").Value = 3
Set-TextValue $ws "B5" "160314"
Set-TextValue $ws "C5" "华夏行业混合(LOF)"
Set-TextValue $ws "D5" "21.89"
Set-TextValue $ws "E5" "92.15"
Set-TextValue $ws "F5" "3.57"
Set-TextValue $ws "G5" "0.7815"
$ws.Range("H5").Value = 10

$ws.Range("A6").Value = 4
Set-TextValue $ws "B6" "519908"
Set-TextValue $ws "C6" "华夏兴华混合A"
Set-TextValue $ws "D6" "9.39"
Set-TextValue $ws "E6" "91.83"
Set-TextValue $ws "F6" "4.62"
Set-TextValue $ws "G6" "0.4338"
$ws.Range("H6").Value = 5

$ws.Range("A7").Value = 5
Set-TextValue $ws "B7" "960004"
Set-TextValue $ws "C7" "华夏兴华混合H"
Set-TextValue $ws "D7" "9.39"
Set-TextValue $ws "E7" "91.83"
Set-TextValue $ws "F7" "4.62"
Set-TextValue $ws "G7" "0.4338"
$ws.Range("H7").Value = 5

$ws.Range("A8").Value = 6
Set-TextValue $ws "B8" "004640"
Set-TextValue $ws "C8" "华夏节能环保股票"
Set-TextValue $ws "D8" "5.45"
Set-TextValue $ws "E8" "90.73"
Set-TextValue $ws "F8" "7.73"
Set-TextValue $ws "G8" "0.4213"
$ws.Range("H8").Value = 2

$ws.Range("A9").Value = 7
Set-TextValue $ws "B9" "001045"
Set-TextValue $ws "C9" "华夏可转债增强债券A"
Set-TextValue $ws "D9" "24.78"
Set-TextValue $ws "E9" "43.12"
Set-TextValue $ws "F9" "1.65"
Set-TextValue $ws "G9" "0.4089"
$ws.Range("H9").Value = 4

$ws.Range("A10").Value = 8
Set-TextValue $ws "B10" "001046"
Set-TextValue $ws "C10" "华夏可转债增强债券I"
Set-TextValue $ws "D10" "24.78"
Set-TextValue $ws "E10" "43.12"
Set-TextValue $ws "F10" "1.65"
Set-TextValue $ws "G10" "0.4089"
$ws.Range("H10").Value = 4

$ws.Range("A11").Value = 9
Set-TextValue $ws "B11" "001113"
Set-TextValue $ws "C11" "南方大数据100指数A"
Set-TextValue $ws "D11" "20.79"
Set-TextValue $ws "E11" "94.23"
Set-TextValue $ws "F11" "1.78"
Set-TextValue $ws "G11" "0.3701"
$ws.Range("H11").Value = 6

$ws.Range("A12").Value = 10
Set-TextValue $ws "B12" "011458"
Set-TextValue $ws "C12" "东方鑫享价值成长一年持有期混合型证券投资基金A"
Set-TextValue $ws "D12" "4.69"
Set-TextValue $ws "E12" "80.98"
Set-TextValue $ws "F12" "2.83"
Set-TextValue $ws "G12" "0.1327"
$ws.Range("H12").Value = 8

$ws.Range("A13").Value = 11
Set-TextValue $ws "B13" "000866"
Set-TextValue $ws "C13" "华宝高端制造股票"
Set-TextValue $ws "D13" "2.25"
Set-TextValue $ws "E13" "91.51"
Set-TextValue $ws "F13" "4.03"
Set-TextValue $ws "G13" "0.0907"
$ws.Range("H13").Value = 1

$ws.Range("A14").Value = 12
Set-TextValue $ws "B14" "012924"
Set-TextValue $ws "C14" "华夏新时代灵活配置混合（QDII）美元现汇"
Set-TextValue $ws "D14" "2.56"
Set-TextValue $ws "E14" "84.71"
Set-TextValue $ws "F14" "3.28"
Set-TextValue $ws "G14" "0.0840"
$ws.Range("H14").Value = 8

$ws.Range("A15").Value = 13
Set-TextValue $ws "B15" "012925"
Set-TextValue $ws "C15" "华夏新时代灵活配置混合（QDII）美元现钞"
Set-TextValue $ws "D15" "2.56"
Set-TextValue $ws "E15" "84.71"
Set-TextValue $ws "F15" "3.28"
Set-TextValue $ws "G15" "0.0840"
$ws.Range("H15").Value = 8

$ws.Range("A16").Value = 14
Set-TextValue $ws "B16" "011459"
Set-TextValue $ws "C16" "东方鑫享价值成长一年持有期混合型证券投资基金C"
Set-TextValue $ws "D16" "1.71"
Set-TextValue $ws "E16" "80.98"
Set-TextValue $ws "F16" "2.83"
Set-TextValue $ws "G16" "0.0484"
$ws.Range("H16").Value = 8

$ws.Range("A17").Value = 15
Set-TextValue $ws "B17" "007770"
Set-TextValue $ws "C17" "同泰开泰混合A"
Set-TextValue $ws "D17" "0.88"
Set-TextValue $ws "E17" "92.35"
Set-TextValue $ws "F17" "3.18"
Set-TextValue $ws "G17" "0.0280"
$ws.Range("H17").Value = 4

$ws.Range("A18").Value = 16
Set-TextValue $ws "B18" "012887"
Set-TextValue $ws "C18" "华夏可转债增强债券C"
Set-TextValue $ws "D18" "0.76"
Set-TextValue $ws "E18" "43.12"
Set-TextValue $ws "F18" "1.65"
Set-TextValue $ws "G18" "0.0125"
$ws.Range("H18").Value = 4

$ws.Range("A19").Value = 17
Set-TextValue $ws "B19" "004284"
Set-TextValue $ws "C19" "华宝新优选一年定期开放灵活配置混合"
Set-TextValue $ws "D19" "0.64"
Set-TextValue $ws "E19" "38.91"
Set-TextValue $ws "F19" "1.69"
Set-TextValue $ws "G19" "0.0108"
$ws.Range("H19").Value = 9

$ws.Range("A20").Value = 18
Set-TextValue $ws "B20" "005128"
Set-TextValue $ws "C20" "华夏永康添福混合"
Set-TextValue $ws "D20" "1.47"
Set-TextValue $ws "E20" "24.52"
Set-TextValue $ws "F20" "0.68"
Set-TextValue $ws "G20" "0.0100"
$ws.Range("H20").Value = 8

$ws.Range("A21").Value = 19
Set-TextValue $ws "B21" "007771"
Set-TextValue $ws "C21" "同泰开泰混合C"
Set-TextValue $ws "D21" "0.18"
Set-TextValue $ws "E21" "92.35"
Set-TextValue $ws "F21" "3.18"
Set-TextValue $ws "G21" "0.0057"
$ws.Range("H21").Value = 4

$ws.Range("A22").Value = 20
Set-TextValue $ws "B22" "004344"
Set-TextValue $ws "C22" "南方大数据100指数C"
Set-TextValue $ws "D22" "0.17"
Set-TextValue $ws "E22" "94.23"
Set-TextValue $ws "F22" "1.78"
Set-TextValue $ws "G22" "0.0030"
$ws.Range("H22").Value = 6


# --- update the "总计" (Total) summary sheet --------------------------------
# Insert a new row 2 for "2022-Q1" above the existing "2021-Q4" row, then
# restyle it to match its sibling data rows (index cell keeps the bordered
# index style; the rest stay unstyled, exactly like rows 3/4 below it).
$total.Rows(2).Insert()

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B3:D3").Copy()
$total.Range("B2:D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 21
$total.Range("D2").Value = 6.63

# the pushed-down rows keep their old running index -- renumber them
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
